# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style used by columns A-H, then fill in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns, copying the style of the existing header row (H1).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# row, I-value, J-value
$data = @(
    @(2,1,5),
    @(3,2,6),
    @(4,7,8),
    @(5,1,3),
    @(6,1,5),
    @(7,8,8),
    @(8,4,6),
    @(9,8,8),
    @(10,3,5),
    @(11,6,8),
    @(12,6,6),
    @(13,2,5),
    @(14,10,11),
    @(15,4,7),
    @(16,5,6),
    @(17,4,6),
    @(18,7,8),
    @(19,7,9),
    @(20,6,8),
    @(21,6,9),
    @(22,2,3),
    @(23,7,9),
    @(24,6,6),
    @(25,5,7),
    @(26,6,9),
    @(27,5,8),
    @(28,8,8),
    @(29,3,8),
    @(30,8,9),
    @(31,7,8),
    @(32,8,9),
    @(33,7,7),
    @(34,4,9),
    @(35,1,4),
    @(36,5,7),
    @(37,7,8),
    @(38,1,3),
    @(39,1,2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
